$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A18").Value = "'90991005"
$ws.Range("B18").Value = "Ar Condicionado Split 9000BTUs Frio 220V Series A1 TCL"
$ws.Range("C18").Value = "1,759,00"

$ws.Range("A19").Value = "'89021443"
$ws.Range("B19").Value = "Depósito para Jardim Manor Resina Cinza 3800 L Keter"
$ws.Range("C19").Value = "3,896,90"

$ws.Range("A20").Value = "'89021443"
$ws.Range("B20").Value = "Depósito para Jardim Manor Resina Cinza 3800 L Keter"
$ws.Range("C20").Value = "3,896,90"
